$wb = $excel.ActiveWorkbook

# Rename sheets: "CS-I2" -> "C-I2", "CS-I3" -> "C-I3"
$wb.Worksheets.Item("CS-I2").Name = "C-I2"
$wb.Worksheets.Item("CS-I3").Name = "C-I3"

# Update selection on the active sheet (C-I3) from F26 to L8
$ws2 = $wb.Worksheets.Item("C-I3")
$ws2.Activate()
$ws2.Range("L8").Select()
